# Add "current URL" validation row to the RegisterNewUser sheet and make
# that sheet the active tab (matches commit "Added current URL validation.
# Documentation is in progress.").

$wb = $excel.ActiveWorkbook

# RegisterNewUser is sheet6.xml / the 6th sheet (0-based index 5), so making
# it the active sheet naturally produces workbookView/@activeTab="5" and
# flips sheetView/@tabSelected from the previously-active "openBrowser"
# sheet onto this one.
$ws = $wb.Worksheets.Item("RegisterNewUser")
$ws.Activate()

# New row 14: label in column A, expected page URL in column B.
$ws.Range("A14").Value = "pageURL"
$ws.Range("B14").Value = "http://newtours.demoaut.com/create_account_success.php"

# Make B14 a clickable hyperlink, like the other URL/email cells on this
# sheet (B4, B12, B13).
$ws.Hyperlinks.Add($ws.Range("B14"), "http://newtours.demoaut.com/create_account_success.php")

# Re-apply the built-in Hyperlink cell style so B14 matches the existing
# hyperlink cells' formatting.
$ws.Range("B14").Style = "Hyperlink"

# Leave B14 selected, as the last-edited cell on the now-active sheet.
[void]$ws.Range("B14").Select()
